$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format (style) from the previous date cell (B8) onto B9
# so it reuses the existing numFmtId=14 (m/d/yyyy) style instead of
# creating a brand-new custom number format.
$ws.Range("B8").Copy($ws.Range("B9"))

# Day 8 row of data
$ws.Range("A9").Value = "Day 8"
$ws.Range("B9").Value = 45810
$ws.Range("C9").Value = "Merge Sorted Arrays"
$ws.Range("D9").Value = "Remove Duplicates from Sorted Array"
$ws.Range("E9").Value = "Remove Duplicates from Sorted Array II"
$ws.Range("F9").Value = "Arrays, Two Pointers"
$ws.Range("G9").Value = "S"
$ws.Range("H9").Value = "YES"

# Move the active selection the same way the source workbook shows
# (cursor sitting just below the newly added data).
$ws.Range("A10").Select()
